# Generate Report for Handoff
# Updates the localization-status workbook so that the zh-cn and de-de
# translations are marked "Ready for handoff" (instead of "In Translation"),
# and the corresponding handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- zh-cn sheet: Status (C2) + Latest Handoff Datetime (E2) ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-24 12:48:05"

# --- de-de sheet: Status (C2) + Latest Handoff Datetime (E2) ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-24 12:48:10"

# --- Overview sheet: per-language status (B2/C2) + Latest Handoff Date (D2) ---
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-24 12:48:10"
